$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 and J1 with same text style as existing headers.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold font, thin border, centered alignment) from an
# existing header cell (H1) onto the two new header cells so they share the
# same style index as the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the I0 / IF data columns for rows 2-65.
$iValues = @(6,7,9,7,7,8,9,9,7,9,9,8,9,9,9,9,9,9,9,9,9,9,9,9,8,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,10,9,9,9,9,9,9,9,9,8,10,9,9,9,9,8,9,9,9,9,9,9,4,5)
$jValues = @(6,7,9,7,7,8,9,9,7,9,9,8,9,9,9,9,9,9,9,10,9,9,9,9,8,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,10,9,9,9,9,9,9,9,9,8,10,9,9,9,9,8,9,9,9,9,9,9,4,5)

for ($k = 0; $k -lt $iValues.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$k]
    $ws.Cells.Item($row, 10).Value = $jValues[$k]
}
